$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory")
$ws.Range("C4").Value = "Car"
$ws.Range("F4").Value = "123123"
